{"js": "// The exercise hint paragraph originally reads (in one run):\n//   \"\u0438\u0437\u043f\u043e\u043b\u0437\u0432\u0430\u0439\u0442\u0435 \u0441\u044a\u0435\u0434\u0438\u043d\u0435\u043d\u0438\u0435 \u043c\u0435\u0436\u0434\u0443 \u0442\u0430\u0431\u043b\u0438\u0446\u0438\u0442\u0435 \"\n// The edit bolds the single word \"\u0441\u044a\u0435\u0434\u0438\u043d\u0435\u043d\u0438\u0435\" inside that sentence,\n// which means splitting the run into three runs:\n//   1) \"\u0438\u0437\u043f\u043e\u043b\u0437\u0432\u0430\u0439\u0442\u0435 \"                (unchanged formatting)\n//   2) \"\u0441\u044a\u0435\u0434\u0438\u043d\u0435\u043d\u0438\u0435\"                  (now bold)\n//   3) \" \u043c\u0435\u0436\u0434\u0443 \u0442\u0430\u0431\u043b\u0438\u0446\u0438\u0442\u0435 \"           (unchanged formatting)\n\nconst searchResults = context.document.body.search(\"\u0441\u044a\u0435\u0434\u0438\u043d\u0435\u043d\u0438\u0435\", { matchCase: true });\nsearchResults.load(\"text\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Could not find the target word '\u0441\u044a\u0435\u0434\u0438\u043d\u0435\u043d\u0438\u0435' to bold.\");\n}\n\n// There is exactly one occurrence of this word in the document (inside the\n// \"\u041f\u043e\u0434\u0441\u043a\u0430\u0437\u043a\u0430\" hint about joining Users/Orders). Bold it in place - Word\n// automatically splits the enclosing run into the three runs described\n// above and keeps the surrounding text/formatting intact.\nconst target = searchResults.items[0];\ntarget.font.bold = true;\n\nawait context.sync();\n", "ps1": "# The exercise hint paragraph originally reads (in one run):\n#   \"\u0438\u0437\u043f\u043e\u043b\u0437\u0432\u0430\u0439\u0442\u0435 \u0441\u044a\u0435\u0434\u0438\u043d\u0435\u043d\u0438\u0435 \u043c\u0435\u0436\u0434\u0443 \u0442\u0430\u0431\u043b\u0438\u0446\u0438\u0442\u0435 \"\n# The edit bolds the single word \"\u0441\u044a\u0435\u0434\u0438\u043d\u0435\u043d\u0438\u0435\" inside that sentence, which\n# means splitting the run into three runs:\n#   1) \"\u0438\u0437\u043f\u043e\u043b\u0437\u0432\u0430\u0439\u0442\u0435 \"                (unchanged formatting)\n#   2) \"\u0441\u044a\u0435\u0434\u0438\u043d\u0435\u043d\u0438\u0435\"                  (now bold)\n#   3) \" \u043c\u0435\u0436\u0434\u0443 \u0442\u0430\u0431\u043b\u0438\u0446\u0438\u0442\u0435 \"           (unchanged formatting)\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"\u0441\u044a\u0435\u0434\u0438\u043d\u0435\u043d\u0438\u0435\"\n$find.MatchCase = $true\n$find.Execute() | Out-Null\n\n# Bold just the matched word - Word splits the surrounding run into the\n# three runs described above and leaves the rest of the sentence as-is.\n$rng.Font.Bold = 1\n"}
